# Case with 380 kV done - update pl_mw.xlsx results (res_line sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.48233129632132
$ws.Range("C2").Value = 0.2482008696125888
$ws.Range("E2").Value = 0.2224249889558898
$ws.Range("F2").Value = 1.730135100785049
$ws.Range("G2").Value = 0.002439348504210898
$ws.Range("I2").Value = 0.7078716785253718
$ws.Range("J2").Value = 0.03852121305155265
$ws.Range("L2").Value = 0.5638261547553753
$ws.Range("O2").Value = 2.668821190617322

# Row 3
$ws.Range("B3").Value = 1.342838457139123
$ws.Range("C3").Value = 0.2336297118149275
$ws.Range("E3").Value = 0.2230072569894173
$ws.Range("F3").Value = 1.73612048323141
$ws.Range("G3").Value = 0.002442154109169425
$ws.Range("I3").Value = 0.7202676486182042
$ws.Range("J3").Value = 0.03710096342302194
$ws.Range("L3").Value = 0.5480674279924074
$ws.Range("O3").Value = 2.698625649367003

# Row 4
$ws.Range("B4").Value = 1.257106423246455
$ws.Range("C4").Value = 0.2246390865130081
$ws.Range("E4").Value = 0.2234407044868014
$ws.Range("F4").Value = 1.740884932627303
$ws.Range("G4").Value = 0.002443968137365772
$ws.Range("I4").Value = 0.7284045246273241
$ws.Range("J4").Value = 0.03623254804905685
$ws.Range("L4").Value = 0.5385549524080204
$ws.Range("O4").Value = 2.719056679568027

# Row 5
$ws.Range("B5").Value = 1.222151400771963
$ws.Range("C5").Value = 0.2209645445754234
$ws.Range("E5").Value = 0.2236364696800575
$ws.Range("F5").Value = 1.743100326541857
$ws.Range("G5").Value = 0.002444730413988642
$ws.Range("I5").Value = 0.7318523656658567
$ws.Range("J5").Value = 0.0358796074998402
$ws.Range("L5").Value = 0.5347199537193035
$ws.Range("O5").Value = 2.727917418158412

# Row 6
$ws.Range("B6").Value = 1.216346095774611
$ws.Range("C6").Value = 0.2203537452764124
$ws.Range("E6").Value = 0.2236701329728756
$ws.Range("F6").Value = 1.743484727579776
$ws.Range("G6").Value = 0.002444858383492399
$ws.Range("I6").Value = 0.7324328437168539
$ws.Range("J6").Value = 0.03582106013499242
$ws.Range("L6").Value = 0.5340856655445521
$ws.Range("O6").Value = 2.729421018310305

# Row 7
$ws.Range("B7").Value = 1.256635079502018
$ws.Range("C7").Value = 0.2245895737154342
$ws.Range("E7").Value = 0.2234432671299551
$ws.Range("F7").Value = 1.740913701558483
$ws.Range("G7").Value = 0.002443978324251651
$ws.Range("I7").Value = 0.7284504892210997
$ws.Range("J7").Value = 0.03622778429424045
$ws.Range("L7").Value = 0.5385030641808441
$ws.Range("O7").Value = 2.71917401391606

# Row 8
$ws.Range("B8").Value = 1.43425270355732
$ws.Range("C8").Value = 0.2431859818369446
$ws.Range("E8").Value = 0.2226100176871739
$ws.Range("F8").Value = 1.73197269428708
$ws.Range("G8").Value = 0.002440296953108719
$ws.Range("I8").Value = 0.7120365848545021
$ws.Range("J8").Value = 0.0380307790792962
$ws.Range("L8").Value = 0.5583588036437419
$ws.Range("O8").Value = 2.678654930298222

# Row 9
$ws.Range("B9").Value = 1.781818881814047
$ws.Range("C9").Value = 0.2792964360062058
$ws.Range("E9").Value = 0.2215769483741141
$ws.Range("F9").Value = 1.723089455463139
$ws.Range("G9").Value = 0.002433799681315419
$ws.Range("I9").Value = 0.6840284774833876
$ws.Range("J9").Value = 0.04159387204904164
$ws.Range("L9").Value = 0.5985813719161968
$ws.Range("O9").Value = 2.616147668264773

# Row 10
$ws.Range("B10").Value = 2.036637107487024
$ws.Range("C10").Value = 0.3055997833247659
$ws.Range("E10").Value = 0.2211823557134203
$ws.Range("F10").Value = 1.721847347271577
$ws.Range("G10").Value = 0.00242946180980208
$ws.Range("I10").Value = 0.6660097676305483
$ws.Range("J10").Value = 0.04422689646723654
$ws.Range("L10").Value = 0.6289048306161078
$ws.Range("O10").Value = 2.580616452829901

# Row 11
$ws.Range("B11").Value = 2.152427113607985
$ws.Range("C11").Value = 0.3175147882221836
$ws.Range("E11").Value = 0.221081573257063
$ws.Range("F11").Value = 1.72243229130612
$ws.Range("G11").Value = 0.002427582064059057
$ws.Range("I11").Value = 0.6583705405436575
$ws.Range("J11").Value = 0.04542773618049267
$ws.Range("L11").Value = 0.6428651272940158
$ws.Range("O11").Value = 2.566721917498711

# Row 12
$ws.Range("B12").Value = 2.19625344309577
$ws.Range("C12").Value = 0.3220192239852793
$ws.Range("E12").Value = 0.2210546940313698
$ws.Range("F12").Value = 1.722819321613812
$ws.Range("G12").Value = 0.002426883637840156
$ws.Range("I12").Value = 0.6555581514671971
$ws.Range("J12").Value = 0.0458828731987424
$ws.Range("L12").Value = 0.6481751381801359
$ws.Range("O12").Value = 2.561787705248349

# Row 13
$ws.Range("B13").Value = 2.186815623630366
$ws.Range("C13").Value = 0.3210494515171831
$ws.Range("E13").Value = 0.2210599816091623
$ws.Range("F13").Value = 1.722728603456162
$ws.Range("G13").Value = 0.002427033461893829
$ws.Range("I13").Value = 0.6561602692055075
$ws.Range("J13").Value = 0.04578483384381116
$ws.Range("L13").Value = 0.6470304900572899
$ws.Range("O13").Value = 2.562835799875501

# Row 14
$ws.Range("B14").Value = 2.15603316338013
$ws.Range("C14").Value = 0.3178855229764679
$ws.Range("E14").Value = 0.2210791359203057
$ws.Range("F14").Value = 1.722460814261112
$ws.Range("G14").Value = 0.002427524336008928
$ws.Range("I14").Value = 0.6581375504720199
$ws.Range("J14").Value = 0.04546517266052064
$ws.Range("L14").Value = 0.6433015149528387
$ws.Range("O14").Value = 2.566309410100388

# Row 15
$ws.Range("B15").Value = 2.137175233335029
$ws.Range("C15").Value = 0.3159465388556555
$ws.Range("E15").Value = 0.2210923370822577
$ws.Range("F15").Value = 1.722318346513688
$ws.Range("G15").Value = 0.002427826752753225
$ws.Range("I15").Value = 0.6593591726287897
$ws.Range("J15").Value = 0.04526942266215883
$ws.Range("L15").Value = 0.6410204673599651
$ws.Range("O15").Value = 2.568479759365459

# Row 16
$ws.Range("B16").Value = 2.029067173963767
$ws.Range("C16").Value = 0.3048200714225118
$ws.Range("E16").Value = 0.2211905231408728
$ws.Range("F16").Value = 1.721832274012129
$ws.Range("G16").Value = 0.002429586533196975
$ws.Range("I16").Value = 0.6665202495943952
$ws.Range("J16").Value = 0.0441484773777745
$ws.Range("L16").Value = 0.6279958052767824
$ws.Range("O16").Value = 2.581570241167896

# Row 17
$ws.Range("B17").Value = 1.962711850461744
$ws.Range("C17").Value = 0.2979812235289785
$ws.Range("E17").Value = 0.2212708942571382
$ws.Range("F17").Value = 1.721828736233974
$ws.Range("G17").Value = 0.002430690023075055
$ws.Range("I17").Value = 0.6710563187628047
$ws.Range("J17").Value = 0.04346157309497301
$ws.Range("L17").Value = 0.6200478849832507
$ws.Range("O17").Value = 2.590182701793964

# Row 18
$ws.Range("B18").Value = 1.92453411048416
$ws.Range("C18").Value = 0.2940429589000075
$ws.Range("E18").Value = 0.2213245331516411
$ws.Range("F18").Value = 1.721934934104624
$ws.Range("G18").Value = 0.002431333532708581
$ws.Range("I18").Value = 0.6737178128662826
$ws.Range("J18").Value = 0.04306677438907514
$ws.Range("L18").Value = 0.6154920955624448
$ws.Range("O18").Value = 2.595349845497594

# Row 19
$ws.Range("B19").Value = 1.911605811620007
$ws.Range("C19").Value = 0.2927087238007573
$ws.Range("E19").Value = 0.2213439683658969
$ws.Range("F19").Value = 1.721989475104621
$ws.Range("G19").Value = 0.002431552929323887
$ws.Range("I19").Value = 0.6746279529399644
$ws.Range("J19").Value = 0.04293315330533432
$ws.Range("L19").Value = 0.6139522802627084
$ws.Range("O19").Value = 2.597135984201714

# Row 20
$ws.Range("B20").Value = 1.969776738435769
$ws.Range("C20").Value = 0.2987097229992628
$ws.Range("E20").Value = 0.2212615717998254
$ws.Range("F20").Value = 1.721817910003921
$ws.Range("G20").Value = 0.002430571643213721
$ws.Range("I20").Value = 0.6705680145536554
$ws.Range("J20").Value = 0.04353466542438866
$ws.Range("L20").Value = 0.6208923378960094
$ws.Range("O20").Value = 2.589243787884101

# Row 21
$ws.Range("B21").Value = 2.165075302658749
$ws.Range("C21").Value = 0.3188150515237567
$ws.Range("E21").Value = 0.2210732038434564
$ws.Range("F21").Value = 1.72253497689421
$ws.Range("G21").Value = 0.002427379791249899
$ws.Range("I21").Value = 0.6575545905201992
$ws.Range("J21").Value = 0.04555905418617101
$ws.Range("L21").Value = 0.644396168759215
$ws.Range("O21").Value = 2.56528023375418

# Row 22
$ws.Range("B22").Value = 2.292591670394131
$ws.Range("C22").Value = 0.3319111098281837
$ws.Range("E22").Value = 0.221015855419882
$ws.Range("F22").Value = 1.723968498353642
$ws.Range("G22").Value = 0.002425371767691309
$ws.Range("I22").Value = 0.6495184295027059
$ws.Range("J22").Value = 0.04688445415450815
$ws.Range("L22").Value = 0.6598943976116232
$ws.Range("O22").Value = 2.551527199669295

# Row 23
$ws.Range("B23").Value = 2.224545788684736
$ws.Range("C23").Value = 0.3249256001348044
$ws.Range("E23").Value = 0.2210404577291953
$ws.Range("F23").Value = 1.723115062861197
$ws.Range("G23").Value = 0.002426436367909957
$ws.Range("I23").Value = 0.6537644989488314
$ws.Range("J23").Value = 0.04617686032779034
$ws.Range("L23").Value = 0.6516102607718892
$ws.Range("O23").Value = 2.558692460278763

# Row 24
$ws.Range("B24").Value = 1.966582793294833
$ws.Range("C24").Value = 0.2983803886799308
$ws.Range("E24").Value = 0.2212657633291428
$ws.Range("F24").Value = 1.721822467418434
$ws.Range("G24").Value = 0.002430625134530711
$ws.Range("I24").Value = 0.6707886097936147
$ws.Range("J24").Value = 0.04350162003140667
$ws.Range("L24").Value = 0.6205105184251636
$ws.Range("O24").Value = 2.58966759896569

# Row 25
$ws.Range("B25").Value = 1.687881447798986
$ws.Range("C25").Value = 0.2695668026249791
$ws.Range("E25").Value = 0.2217922881850643
$ws.Range("F25").Value = 1.724565405948724
$ws.Range("G25").Value = 0.002435480539961658
$ws.Range("I25").Value = 0.6911568394611507
$ws.Range("J25").Value = 0.04062718116349373
$ws.Range("L25").Value = 0.5875636996242264
$ws.Range("O25").Value = 2.631236718433286

